$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95, shifting existing rows 95-99 down to 96-100.
$ws.Rows.Item(95).Insert()

# Populate the newly inserted row 95 with the latest weekly price record.
$ws.Range("A95").Value = 5
$ws.Range("B95").Value = "Macroferia Regional de Talca"
$ws.Range("C95").Value = "Maule"
$ws.Range("D95").Value = 44568
$ws.Range("E95").Value = 7
$ws.Range("F95").Value = "Fruta"
$ws.Range("G95").Value = 100108
$ws.Range("H95").Value = "Tropicales y subtropicales"
$ws.Range("I95").Value = 100108002
$ws.Range("J95").Value = "Mango"
$ws.Range("K95").Value = "Sin especificar"
$ws.Range("L95").Value = "Primera"
$ws.Range("M95").Value = 420
$ws.Range("N95").Value = 6000
$ws.Range("O95").Value = 6000
$ws.Range("P95").Value = 6000
$ws.Range("Q95").Value = "$/bandeja 4 kilos"
$ws.Range("R95").Value = "Perú"
$ws.Range("S95").Value = 1500
$ws.Range("T95").Value = 4

# Make sure the new row's date cell keeps the same date style used by the
# rest of column D (it is already inherited from the Insert, but set it
# explicitly to be safe).
$ws.Range("D95").NumberFormat = $ws.Range("D96").NumberFormat()
